# Fixed data, added availability

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TemperatureRelay")

# --- Fix pricing data on the TemperatureRelay sheet ---
$ws.Range("B2").Value = 190
$ws.Range("B3").Value = 290
$ws.Range("B4").Value = 308
$ws.Range("B5").Value = 650
$ws.Range("B6").Value = 305
$ws.Range("B7").Value = 290
$ws.Range("B8").Value = 426
$ws.Range("B9").Value = 487
$ws.Range("B10").Value = 457
$ws.Range("B11").Value = 586

# Move the active selection on the first sheet
[void]$ws.Range("D14").Select()

# --- Add a new "Metadata" sheet with availability info, placed right after TemperatureRelay ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "Metadata"

$meta.Range("A1").Value = "Энергохит"
$meta.Range("B1").Value = "24.07.2013"
$meta.Range("C1").Value = "0.08.2012"
$meta.Range("D1").Value = "Price update"

[void]$meta.Range("D3").Select()

# Restore TemperatureRelay as the active/selected sheet
[void]$ws.Activate()
